# "more and more footprints"
# Change every "m" marker in column D (the "Auf Lager" / stock-status column)
# to "n" on the HybridChargeController sheet, and move the active selection
# to A13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("HybridChargeController")

$rows = @(9, 10, 11, 13, 15, 16, 17, 18, 19, 21, 29, 30, 32)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq "m") {
        $cell.Value = "n"
    }
}

$ws.Range("A13").Select()
